$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '69.967.93'
$ws.Cells.Item(2, 5).Value = '  -1.38%  '
$ws.Cells.Item(3, 4).Value = '3.753.92'
$ws.Cells.Item(3, 5).Value = '  +2.45%  '
$ws.Cells.Item(4, 5).Value = '  +0.13%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '622.66'
$cell.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.60%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '180.00'
$cell.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.78%  '
$ws.Cells.Item(7, 4).Value = '3.752.70'
$ws.Cells.Item(7, 5).Value = '  +2.50%  '
$ws.Cells.Item(8, 5).Value = '  +0.04%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.534'
$cell.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -1.62%  '
$ws.Cells.Item(10, 5).Value = '  +2.89%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.32'
$cell.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -5.46%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.488'
$cell.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -3.28%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '41.30'
$cell.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +2.21%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0000260'
$cell.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +2.13%  '
$ws.Cells.Item(15, 4).Value = '4.381.34'
$ws.Cells.Item(15, 5).Value = '  +2.52%  '
$ws.Cells.Item(16, 4).Value = '3.764.43'
$ws.Cells.Item(16, 5).Value = '  +3.32%  '
$ws.Cells.Item(17, 4).Value = '70.038.90'
$ws.Cells.Item(17, 5).Value = '  -1.29%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.123'
$cell.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -0.48%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.61'
$cell.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +0.54%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '16.76'
$cell.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -0.86%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '506.49'
$cell.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -2.88%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '9.43'
$cell.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +1.36%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.728'
$cell.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -2.22%  '
$ws.Cells.Item(24, 5).Value = '  -1.24%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '87.00'
$cell.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -1.99%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '13.14'
$cell.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -2.68%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '11.18'
$cell.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.52%  '
$ws.Cells.Item(28, 5).Value = '  +25.92%  '
$ws.Cells.Item(29, 5).Value = '  +0.13%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.49'
$cell.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -2.11%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.95'
$cell.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +1.61%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.91'
$cell.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -3.42%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '31.34'
$cell.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -1.03%  '
$ws.Cells.Item(34, 5).Value = '  -0.11%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +0.17%  '
$ws.Cells.Item(36, 5).Value = '  +4.17%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.22'
$cell.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +1.38%  '
$ws.Cells.Item(38, 5).Value = '  -4.88%  '
$ws.Cells.Item(39, 5).Value = '  +1.32%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.11'
$cell.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -3.84%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '50.37'
$cell.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -1.93%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '45.02'
$cell.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -2.74%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '423.66'
$cell.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -1.04%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '8.75'
$cell.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.80%  '
$ws.Cells.Item(45, 5).Value = '  +2.17%  '
$ws.Cells.Item(46, 4).Value = '3.004.99'
$ws.Cells.Item(46, 5).Value = '  -3.71%  '
$ws.Cells.Item(47, 5).Value = '  -1.24%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '27.34'
$cell.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -3.70%  '
$ws.Cells.Item(49, 5).Value = '  -0.04%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '138.16'
$cell.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -1.88%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.52'
$cell.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +1.21%  '
